$d = $word.ActiveDocument

# 1. Remove the first "_GoBack" bookmark (the one in the skills table cell,
#    right after "MS SQL, JavaScript ").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Update the "Try to Solve ..." sentence wording.
$d.Content.Find.Execute("Try to Solve (SQL & Python) problem at ", $false, $false, $false, $false, $false, $true, 1, $false, "Solving code (SQL & Python) challenges on ", 2) | Out-Null

# 3. Reword the "Gold Level" sentence.
$d.Content.Find.Execute(". 58/58 challenges solved (Gold Level) in SQL as well as ", $false, $false, $false, $false, $false, $true, 1, $false, ". Achieved Gold Level (58/58 challenges solved) in SQL as well as ", 2) | Out-Null

# 4. Merge the two paragraphs ("... as well as" / "Silver Level in Python.")
#    into a single paragraph by replacing the paragraph mark between them
#    with nothing.
$d.Content.Find.Execute("as well as " + [char]13 + "Silver", $false, $false, $false, $false, $false, $true, 1, $false, "as well as Silver", 2) | Out-Null

# 5. Re-add the "_GoBack" bookmark at the very end of the now-merged
#    paragraph (right after "Silver Level in Python.").
$target = $d.Content.Find.Execute("Silver Level in Python.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r = $d.Content
$r.Find.Execute("Silver Level in Python.") | Out-Null
$endRange = $d.Range($r.Find.Parent.End, $r.Find.Parent.End)
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null
